# edit vessel size SE in (Tufto & Rofstad, 1999)
#
# On the "Vessel size (tumor)" sheet, rows 6-8 hold the Tufto & Rofstad,
# 1999 (D-12 / R-18 / U-25 cell) entries; update their SE (column C)
# values. Also replay the sheet/cell selection state (which sheet/cell
# was last active) left behind by the edit.

$wb = $excel.ActiveWorkbook

# --- Update the SE values for the Tufto & Rofstad, 1999 rows ---
$wsTumor = $wb.Worksheets.Item("Vessel size (tumor)")
$wsTumor.Range("C6").Value = 8.32
$wsTumor.Range("C7").Value = 4.23
$wsTumor.Range("C8").Value = 5.51

# --- Replay the selection left on each visited sheet ---
$wsAdiposeSize = $wb.Worksheets.Item("Vessel size (adipose)")
$wsAdiposeSize.Activate() | Out-Null
$wsAdiposeSize.Range("C17").Select() | Out-Null

$wsAdiposeDensity = $wb.Worksheets.Item("Vessel density (adipose)")
$wsAdiposeDensity.Activate() | Out-Null
$wsAdiposeDensity.Range("A6").Select() | Out-Null

$wsTumorDensity = $wb.Worksheets.Item("Vessel density (tumor)")
$wsTumorDensity.Activate() | Out-Null
$wsTumorDensity.Range("B5").Select() | Out-Null

# Finish on the "Vessel size (tumor)" sheet, which is the active tab
# in the saved workbook, with the last-edited cell selected.
$wsTumor.Activate() | Out-Null
$wsTumor.Range("C8").Select() | Out-Null
